$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 63
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 14
$ws.Range("E2").Value = 94
$ws.Range("F2").Value = 6
